$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "NA" under the duplicate_image_filename column (E) for the main
# stimuli table rows (practice rows 2-5, generic rows 6-13, and the
# unique_video/unique_audio rows 14-21).
$ws.Range("E2:E21").Value = "NA"
